# Updated cryptos list (refresh of price / volume(1h) columns, plus a
# Kaspa <-> FirstDigitalUSD row swap around rank ~32/33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while forcing text storage (so values
# like "597.41" or "0.999" keep their exact textual form instead of being
# auto-coerced to a number). The leading apostrophe is Excel's classic
# "treat as text" marker; we restore the cell's original style right after
# so we don't leave a stray quote-prefix style behind.
function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.Value = "'" + $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "68.257.49"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.647.90"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "597.41"
$ws.Range("E5").Value = "  -0.42%  "
Set-TextValue "D6" "156.70"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +0.68%  "
Set-TextValue "D13" "28.03"
$ws.Range("E13").Value = "  +0.70%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "3.129.67"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "68.143.48"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "2.649.07"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("E18").Value = "  -0.16%  "
Set-TextValue "D19" "363.43"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("E23").Value = "  -2.10%  "
Set-TextValue "D24" "75.18"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "2.779.88"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("E28").Value = "  -0.36%  "
Set-TextValue "D29" "0.999"
$ws.Range("E29").Value = "  -0.03%  "
Set-TextValue "D30" "559.92"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D34" "0.128"
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D35" "0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.58%  "
Set-TextValue "D37" "161.85"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  +1.24%  "
$ws.Range("E40").Value = "  -2.73%  "
Set-TextValue "D41" "5.33"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("E42").Value = "  +4.22%  "
Set-TextValue "D43" "17.80"
$ws.Range("E43").Value = "  +0.32%  "
Set-TextValue "D44" "2.61"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E45").Value = "  +0.06%  "
Set-TextValue "D46" "158.96"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  -0.34%  "
Set-TextValue "D48" "22.11"
$ws.Range("E48").Value = "  +0.79%  "
Set-TextValue "D49" "1.69"
$ws.Range("E49").Value = "  -1.34%  "
Set-TextValue "D50" "0.0781"
$ws.Range("E50").Value = "  +0.37%  "
Set-TextValue "D51" "0.615"
$ws.Range("E51").Value = "  -0.15%  "
